# Generate Report for Handback
# Updates the "generated date/time" timestamp strings written into the
# handback-status report after a (re-)run of the report generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 13:10:56"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 13:10:51"
$wsZhCn.Range("K2").Value = "2016-09-03 13:11:16"

# --- de-de sheet --------------------------------------------------------
# de-de!H2 originally shared the same text as Overview!G2
# ("2016-09-03 13:10:15"), so it moves to the same new value ("...13:10:56").
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 13:10:56"
$wsDeDe.Range("K2").Value = "2016-09-03 13:11:23"
